$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.084.34'
$ws.Range('E2').Value = '  +2.16%  '
$ws.Range('D3').Value = '2.417.28'
$ws.Range('E3').Value = '  +2.88%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'555.62"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').Value = "'142.71"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.72%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +1.92%  '
$ws.Range('D9').Value = '2.415.52'
$ws.Range('E9').Value = '  +2.86%  '
$ws.Range('E10').Value = '  +3.84%  '
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('E12').Value = '  +1.56%  '
$ws.Range('E13').Value = '  +1.72%  '
$ws.Range('D14').Value = "'26.17"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +6.21%  '
$ws.Range('D15').Value = "'0.0000173"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +7.90%  '
$ws.Range('D16').Value = '2.856.99'
$ws.Range('E16').Value = '  +3.00%  '
$ws.Range('D17').Value = '62.069.20'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').Value = '2.416.27'
$ws.Range('E18').Value = '  +2.72%  '
$ws.Range('D19').Value = "'11.05"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.01%  '
$ws.Range('E20').Value = '  +1.70%  '
$ws.Range('D21').Value = "'323.56"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.21%  '
$ws.Range('D22').Value = "'6.69"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.03%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = "'64.95"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('D25').Value = "'1.76"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +5.87%  '
$ws.Range('E26').Value = '  +8.04%  '
$ws.Range('D27').Value = "'576.92"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +17.15%  '
$ws.Range('E28').Value = '  +3.09%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').Value = "'8.37"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.77%  '
$ws.Range('D31').Value = '0.0₃0929'
$ws.Range('E31').Value = '  +7.55%  '
$ws.Range('D32').Value = "'1.46"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +6.18%  '
$ws.Range('E33').Value = '  +1.39%  '
$ws.Range('E34').Value = '  +3.72%  '
$ws.Range('E35').Value = '  +3.02%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = "'4.82"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.27%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value = "'5.65"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +8.10%  '
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = "'1.87"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.08%  '
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').Value = "'18.74"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('D42').Value = "'148.15"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.09%  '
$ws.Range('D44').Value = "'41.73"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'150.96"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.21%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').Value = "'2.29"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +12.01%  '
$ws.Range('E47').Value = '  +1.84%  '
$ws.Range('D48').Value = "'0.0543"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.39%  '
$ws.Range('D49').Value = "'20.32"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +6.63%  '
$ws.Range('D50').Value = "'0.587"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.45%  '
$ws.Range('E51').Value = '  +1.97%  '
